$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry below is one data row (rows 2-51 on the sheet): the new "Price"
# (column D) text - or $null if that row's price did not change in this
# update - and the new "Volume(1h)" (column E) text, already padded with the
# same two leading/trailing spaces used throughout the sheet.
$updates = @(
    @{ Row = 2; D = '27.902.15'; E = '  +3.20%  ' }
    @{ Row = 3; D = '1.724.36'; E = '  +2.82%  ' }
    @{ Row = 4; D = '1.00'; E = '  -0.27%  ' }
    @{ Row = 5; D = '217.20'; E = '  +1.04%  ' }
    @{ Row = 6; D = '0.522'; E = '  +1.07%  ' }
    @{ Row = 7; D = '1.00'; E = '  -0.15%  ' }
    @{ Row = 8; D = '23.71'; E = '  +11.30%  ' }
    @{ Row = 9; D = '0.264'; E = '  +3.45%  ' }
    @{ Row = 10; D = $null; E = '  +1.34%  ' }
    @{ Row = 11; D = '0.0898'; E = '  +1.03%  ' }
    @{ Row = 12; D = '1.969.60'; E = '  +3.00%  ' }
    @{ Row = 13; D = '1.719.92'; E = '  +2.39%  ' }
    @{ Row = 14; D = $null; E = '  +3.03%  ' }
    @{ Row = 15; D = '0.566'; E = '  +5.64%  ' }
    @{ Row = 16; D = '67.86'; E = '  +2.49%  ' }
    @{ Row = 17; D = '27.878.76'; E = '  +3.22%  ' }
    @{ Row = 18; D = '241.19'; E = '  +2.37%  ' }
    @{ Row = 19; D = '7.96'; E = '  -3.48%  ' }
    @{ Row = 20; D = $null; E = '  +1.71%  ' }
    @{ Row = 21; D = '1.00'; E = '  -0.27%  ' }
    @{ Row = 22; D = '4.62'; E = '  +3.14%  ' }
    @{ Row = 23; D = '9.71'; E = '  +4.68%  ' }
    @{ Row = 24; D = $null; E = '  +0.26%  ' }
    @{ Row = 25; D = '148.76'; E = '  +1.67%  ' }
    @{ Row = 26; D = '7.52'; E = '  +3.69%  ' }
    @{ Row = 27; D = '16.63'; E = '  +1.57%  ' }
    @{ Row = 28; D = $null; E = '  +1.11%  ' }
    @{ Row = 29; D = '1.00'; E = '  -0.03%  ' }
    @{ Row = 30; D = '0.0508'; E = '  +1.89%  ' }
    @{ Row = 31; D = $null; E = '  +1.31%  ' }
    @{ Row = 32; D = '3.44'; E = '  +2.08%  ' }
    @{ Row = 33; D = '3.29'; E = '  +3.64%  ' }
    @{ Row = 34; D = '1.468.07'; E = '  -4.55%  ' }
    @{ Row = 35; D = '1.67'; E = '  -1.84%  ' }
    @{ Row = 36; D = '0.961'; E = '  +5.29%  ' }
    @{ Row = 37; D = '0.610'; E = '  +3.28%  ' }
    @{ Row = 38; D = $null; E = '  +0.73%  ' }
    @{ Row = 39; D = $null; E = '  -0.21%  ' }
    @{ Row = 40; D = $null; E = '  -1.29%  ' }
    @{ Row = 41; D = '71.39'; E = '  +5.56%  ' }
    @{ Row = 42; D = '5.85'; E = '  +6.26%  ' }
    @{ Row = 43; D = '1.00'; E = '  -0.25%  ' }
    @{ Row = 44; D = '1.873.80'; E = '  +3.01%  ' }
    @{ Row = 45; D = '2.28'; E = '  +0.95%  ' }
    @{ Row = 46; D = '0.789'; E = '  +1.07%  ' }
    @{ Row = 47; D = $null; E = '  +9.05%  ' }
    @{ Row = 48; D = '91.82'; E = '  +1.38%  ' }
    @{ Row = 49; D = $null; E = '  +4.86%  ' }
    @{ Row = 50; D = '8.24'; E = '  +2.26%  ' }
    @{ Row = 51; D = '0.105'; E = '  +1.52%  ' }
)

# Matches a plain signed decimal number such as "1.00" or "217.20" - values
# Excel would otherwise silently re-parse as a number (losing e.g. a trailing
# zero). Values like "27.902.15" have multiple dots, so Excel already treats
# them as plain text and need no special handling.
$numericPattern = '^-?[0-9]+(\.[0-9]+)?$'

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.D -match $numericPattern) {
            # Prefix with a single quote, exactly like typing '1.00 into Excel,
            # to force the literal text to be kept instead of being converted
            # into the number 1.
            $dCell.Formula = "'" + $u.D
        } else {
            $dCell.Value = $u.D
        }
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
